# Add new rows 31-35 to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the existing date-formatted style from C2:D2 down into the new rows'
# date columns (C31:D35) so the new cells reuse the workbook's existing
# date style instead of minting a new numFmt/style entry.
$ws.Range("C2:D2").Copy($ws.Range("C31:D31"))
$ws.Range("C2:D2").Copy($ws.Range("C32:D32"))
$ws.Range("C2:D2").Copy($ws.Range("C33:D33"))
$ws.Range("C2:D2").Copy($ws.Range("C34:D34"))
$ws.Range("C2:D2").Copy($ws.Range("C35:D35"))

# Row 31
$ws.Range("A31").Value = 9876543
$ws.Range("B31").Value = "Performance Tune"
$ws.Range("C31").Value = [DateTime]::FromOADate(42288)
$ws.Range("D31").Value = [DateTime]::FromOADate(42288)
$ws.Range("E31").Value = "RG"
$ws.Range("F31").Value = 7654326

# Row 32
$ws.Range("A32").Value = 7654
$ws.Range("B32").Value = "Performance Tune"
$ws.Range("C32").Value = [DateTime]::FromOADate(42288)
$ws.Range("D32").Value = [DateTime]::FromOADate(42288)
$ws.Range("E32").Value = "RG"

# Row 33
$ws.Range("A33").Value = 456789
$ws.Range("B33").Value = "Performance Tune"
$ws.Range("C33").Value = [DateTime]::FromOADate(42045)
$ws.Range("D33").Value = [DateTime]::FromOADate(42288)
$ws.Range("E33").Value = "RG"
$ws.Range("F33").Value = 9876543

# Row 34
$ws.Range("A34").Value = 4567890
$ws.Range("B34").Value = "Performance Tune"
$ws.Range("C34").Value = [DateTime]::FromOADate(42015)
$ws.Range("D34").Value = [DateTime]::FromOADate(42015)
$ws.Range("E34").Value = "DM"
$ws.Range("F34").Value = 98765

# Row 35
$ws.Range("A35").Value = "4567890-"
$ws.Range("B35").Value = "Accident Calibration"
$ws.Range("C35").Value = [DateTime]::FromOADate(42401)
$ws.Range("D35").Value = [DateTime]::FromOADate(42401)
$ws.Range("E35").Value = "RG"
$ws.Range("F35").Value = 567890
